$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 now describes the new "CloneScene" entry.
$ws.Range("A2").Value = "../../NFDataCfg/Ini/NFZoneServer/Scene/CloneScene/"
$ws.Range("B2").Value = "3"
$ws.Range("F2").Value = "clone"

# Row 3 / Row 4 SceneName values were renamed from "Stage001" to "newscene".
$ws.Range("F3").Value = "newscene"
$ws.Range("F4").Value = "newscene"

# Move/restore the sheet's active selection to H8.
$ws.Range("H8").Select()
